# Update "想去人数" (F column) counts on the "展览", "演出" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 621
$ws1.Range("F5").Value = 162
$ws1.Range("F6").Value = 9399
$ws1.Range("F9").Value = 1199
$ws1.Range("F10").Value = 1138
$ws1.Range("F11").Value = 148
$ws1.Range("F12").Value = 96
$ws1.Range("F13").Value = 17
$ws1.Range("F15").Value = 416
$ws1.Range("F17").Value = 251
$ws1.Range("F18").Value = 1278

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 11

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 11
$ws4.Range("F5").Value = 621
$ws4.Range("F6").Value = 162
$ws4.Range("F7").Value = 9399
$ws4.Range("F10").Value = 1199
$ws4.Range("F11").Value = 1138
$ws4.Range("F12").Value = 148
$ws4.Range("F13").Value = 96
$ws4.Range("F14").Value = 17
$ws4.Range("F16").Value = 416
$ws4.Range("F18").Value = 251
$ws4.Range("F19").Value = 1278
